# "Averaged Data Table created"
#
# Fills the (previously empty) data columns of the "Data" / Table1 table
# with placeholder "N/A" values for every consensus-mechanism row, carries
# over two known numeric data points (TPS = 1763, Nakamoto Coefficient = 18)
# for "Proof of History" and "Proof of History with Proof of Stake", and
# applies per-column number formats (0.00 for the numeric columns, 0% for
# the percentage column, and Text for the two free-text columns) matching
# how the columns are used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric-ish columns (TPS, Energy Use per Transaction, Nakamoto Coefficient) ---
$ws.Range("B2:D10").Value2 = "N/A"
$ws.Range("B2:D10").NumberFormat = "0.00"

# --- Percentage column (% of nodes required to take over network) ---
$ws.Range("E2:E10").Value2 = "N/A"
$ws.Range("E2:E10").NumberFormat = "0%"

# --- Free-text columns (Strengths, Weaknesses) ---
$ws.Range("F2:G10").Value2 = "N/A"
$ws.Range("F2:G10").NumberFormat = "@"

# --- Known data points ---
# Proof of History (row 5) and Proof of History with Proof of Stake (row 7)
$ws.Range("B5").Value2 = 1763
$ws.Range("D5").Value2 = 18
$ws.Range("B7").Value2 = 1763
$ws.Range("D7").Value2 = 18

# --- Restore the selection to where the editor ended up ---
$ws.Range("J10").Select()
